$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing extr1..extr8 rows (currently rows 8-15) down by two
# rows, to rows 10-17, to make room for two new "line7"/"line8" rows.
# Copy bottom-up so we never clobber a row before it has been read.
for ($r = 15; $r -ge 8; $r--) {
    $newR = $r + 2
    $ws.Range("A$r" + ":E$r").Copy($ws.Range("A$newR"))
}

# --- Step 2: write the two new rows (line7, line8) into the now-vacated rows 8-9
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Step 3: fix up the "in_service" flag (and A index) for the shifted rows
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("E12").Value = $false

$ws.Range("A13").Value = 11
$ws.Range("E13").Value = $true

$ws.Range("A14").Value = 12
$ws.Range("E14").Value = $true

$ws.Range("A15").Value = 13
$ws.Range("E15").Value = $true

$ws.Range("A16").Value = 14
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("E17").Value = $true
